# Word COM-interop script applying the target edit to the document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Project " + "Report on" -> "Project Report on"
#    (two runs merge into a single run with no xml:space="preserve")
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Project Report on", $false, $false, $false, $false, $false, $true, 1, $false, "Project Report on", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Time of " + "Year Segmentation" -> "Time of Year Segmentation"
#    (keeps the Strong run style, merges the two Strong runs)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Time of Year Segmentation", $false, $false, $false, $false, $false, $true, 1, $false, "Time of Year Segmentation", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2b. Merge the many small runs describing the segmentation periods into one.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Segment the ride into periods (e.g., summer, winter, evening, fall, spring).", $false, $false, $false, $false, $false, $true, 1, $false, "Segment the ride into periods (e.g., summer, winter, evening, fall, spring).", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Merge the two "Distribution analysis ..." runs into a single run.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Distribution analysis refers to the process of examining the statistical distribution of data within a dataset. Extensive usage of histograms and graphs for distribution analysis.", $false, $false, $false, $false, $false, $true, 1, $false, "Distribution analysis refers to the process of examining the statistical distribution of data within a dataset. Extensive usage of histograms and graphs for distribution analysis.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Insert the new "Clustering / Dimensionality Reduction / Feature
#    Selection / ... / Model Selection and Training" subsections right
#    after the "Data Analysis and Modelling" heading.
#    InsertXML on a collapsed range replaces the whole enclosing paragraph
#    in this runtime, so we rebuild that paragraph (identical to the
#    original) followed by the five new paragraphs, then replace the
#    original paragraph's range with that combined fragment.
# ---------------------------------------------------------------------------
$headingPara = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -eq "Data Analysis and Modelling`r") {
        $headingPara = $cand
        break
    }
}
if ($headingPara -ne $null) {
    $rng = $headingPara.Range
    $inner = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="10" w:name="_Toc171364809"/><w:r><w:t>Data Analysis and Modelling</w:t></w:r><w:bookmarkEnd w:id="10"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Clustering</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Dimensionality Reduction</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Feature Selection</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Feature selection is a critical process in data analysis and modeling, aimed at identifying the most relevant features that contribute to the predictive power of a model. Initially, a correlation analysis is conducted to identify and remove highly correlated features, reducing redundancy in the dataset. Subsequently, recursive feature elimination (RFE) is applied to further refine the feature set. RFE works by recursively fitting a model and removing the least important features, as determined by the model''s performance. This iterative process continues until the optimal subset of features is identified, balancing model accuracy and complexity. To handle large datasets efficiently and manage memory usage, the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>data.table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package is used. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>data.table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> provides optimized and fast data manipulation capabilities, which are particularly useful when dealing with large volumes of data. By selecting a subset of the most significant features and using </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>data.table</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for efficient data processing, we enhance the model''s interpretability and efficiency, ultimately improving its predictive performance.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Model Selection and Training</w:t></w:r></w:p>'
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $inner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 5. Remove <w:lastRenderedPageBreak/> from the "Conclusion" heading run,
#    rebuilding that paragraph without the stray rendering-break marker.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$conclusionPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -eq "Conclusion`r") {
        $conclusionPara = $cand
        break
    }
}
if ($conclusionPara -ne $null) {
    $rng = $conclusionPara.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="12" w:name="_Toc171364811"/><w:r><w:t>Conclusion</w:t></w:r><w:bookmarkEnd w:id="12"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 6. Add <w:lastRenderedPageBreak/> right before "Sustainability Impact:"
#    (it moved here from the Conclusion heading), rebuilding that paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$sustainPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "Sustainability Impact:*") {
        $sustainPara = $cand
        break
    }
}
if ($sustainPara -ne $null) {
    $rng = $sustainPara.Range
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Sustainability Impact:</w:t></w:r><w:r><w:t xml:space="preserve"> Evaluate the environmental benefits of the bike-sharing program and suggest improvements for increasing its positive impact on urban sustainability.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

Write-Output "edit complete"
